$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("S2").Value = 1.5
$ws.Range("T2").Value = 2.37

# Row 3
$ws.Range("M3").Value = 1.1
$ws.Range("N3").Value = 7
$ws.Range("O3").Value = 1.5
$ws.Range("P3").Value = 2.5
$ws.Range("S3").Value = 1.5
$ws.Range("T3").Value = 2.37

# Row 4
$ws.Range("S4").Value = 1.54

# Row 5
$ws.Range("I5").Value = 3.7
$ws.Range("J5").Value = 2.75
$ws.Range("L5").Value = 4.5
$ws.Range("N5").Value = 7.5
$ws.Range("S5").Value = 1.47
$ws.Range("X5").Value = 8.5
$ws.Range("Z5").Value = 17
$ws.Range("AX5").Value = 23
$ws.Range("BA5").Value = 126
